# "contents moved: dialog/ -> /index.html"
#
# The sitemap listed three pages (dialogue, about, articles) that have now
# been moved/merged away; their rows are removed from the sitemap sheet,
# leaving only the "home" ( / ) entry and the trailing EndOfData marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sitemap rows for /dialogue/, /about/ and /articles/{*}
# (rows 10-12). This also shifts the trailing "EndOfData" marker row
# up from row 15 to row 12, and shrinks the used range accordingly.
$ws.Rows("10:12").Delete()

# Cosmetic touches that came along with the resave in Excel:
# gridlines visible again, the title row a bit taller, and the cursor
# left on F19.
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Rows("3").RowHeight = 38
$ws.Range("F19").Select() | Out-Null
